$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# New list of "Rescatables" (students pending retake) after inserting a new
# student (MUÑOZ REYES ERWIN ISRAEL) at the top of the list and moving
# TORRES GUTIERREZ JESUS ENRIQUE to the bottom of the list.
$rows = @(
    @{ A = 24330051920393; B = "MUÑOZ";    C = "REYES";      D = "ERWIN ISRAEL";     E = "Pensamiento matemático II"; F = "2APV"; G = 4 },
    @{ A = 24330051920274; B = "CLEMENTE"; C = "JUAREZ";      D = "BRYAN";            E = "Pensamiento matemático II"; F = "2APV"; G = 3 },
    @{ A = 24330051920255; B = "LOPEZ";    C = "ROSAS";       D = "ERNESTO";          E = "Pensamiento matemático II"; F = "2APV"; G = 3 },
    @{ A = 23330051920312; B = "VERA";     C = "PEREZ";       D = "ALEYDA MONSERRAT"; E = "Pensamiento matemático II"; F = "2ASV"; G = 3 },
    @{ A = 22330051920007; B = "CARRERA";  C = "GARCIA";      D = "ANA KAREN";        E = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"; F = "4BEM"; G = 2 },
    @{ A = 24330051920369; B = "TORRES";   C = "GUTIERREZ";   D = "JESUS ENRIQUE";    E = "Pensamiento matemático II"; F = "2APV"; G = 1 }
)

# Fill column-by-column (B, then C, then D) to match the original
# data-entry order, then the remaining columns.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $rows[$i].B
}
for ($i = 0; $i -lt $rows.Count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $rows[$i].C
}
for ($i = 0; $i -lt $rows.Count; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $rows[$i].D
}
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $rows[$i].A
    $ws.Cells.Item($r, 5).Value = $rows[$i].E
    $ws.Cells.Item($r, 6).Value = $rows[$i].F
    $ws.Cells.Item($r, 7).Value = $rows[$i].G
}
